$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 907
$ws.Range("I2").Value = 338.75
$ws.Range("K2").Value = 338.75
$ws.Range("M2").Value = -225.75

$ws.Range("H19").Value = 2480.7812
$ws.Range("J19").Value = 2757.158
$ws.Range("L19").Value = 2757.158
$ws.Range("N19").Value = -3107.158

$ws.Range("H29").Value = 8872.375
$ws.Range("I29").Value = 4499
$ws.Range("K29").Value = 13497
$ws.Range("M29").Value = -13216

$ws.Range("H106").Value = 7565
$ws.Range("I106").Value = 7565
$ws.Range("K106").Value = 7565
$ws.Range("M106").Value = -6934

$ws.Range("H116").Value = 3592.25
$ws.Range("I116").Value = 2752.0908
$ws.Range("K116").Value = 2752.0908
$ws.Range("M116").Value = 689.9092000000001

$ws.Range("H125").Value = 1947.5
$ws.Range("J125").Value = 1000
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -13920

$ws.Range("H137").Value = 4055.2
$ws.Range("I137").Value = 4003.2
$ws.Range("K137").Value = 12009.6
$ws.Range("M137").Value = -9459.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7109.769
$ws.Range("I61").Value = 6315.1113
$ws.Range("K61").Value = 6315.1113
$ws.Range("M61").Value = -6103.1113

$ws.Range("H122").Value = 7638.516
$ws.Range("I122").Value = 6878.5
$ws.Range("J122").Value = 10244.286
$ws.Range("K122").Value = 20635.5
$ws.Range("L122").Value = 30732.858
$ws.Range("M122").Value = -18185.5
$ws.Range("N122").Value = -35632.858

$ws.Range("H136").Value = 7109.769
$ws.Range("I136").Value = 6315.1113
$ws.Range("K136").Value = 18945.3339
$ws.Range("M136").Value = -16395.3339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3218.054
$ws.Range("I134").Value = 2973
$ws.Range("K134").Value = 8919
$ws.Range("M134").Value = -6384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7296.7915
$ws.Range("I99").Value = 6283.3125
$ws.Range("J99").Value = 9323.75
$ws.Range("K99").Value = 6283.3125
$ws.Range("L99").Value = 9323.75
$ws.Range("M99").Value = -4785.3125
$ws.Range("N99").Value = -12319.75

$ws.Range("H107").Value = 23810810
$ws.Range("I107").Value = 33333958
$ws.Range("J107").Value = 2941.3333
$ws.Range("K107").Value = 33333958
$ws.Range("L107").Value = 2941.3333
$ws.Range("M107").Value = -33332038
$ws.Range("N107").Value = -6781.3333

$ws.Range("H122").Value = 4878.2
$ws.Range("I122").Value = 326.42856
$ws.Range("J122").Value = 15499
$ws.Range("K122").Value = 979.28568
$ws.Range("L122").Value = 46497
$ws.Range("M122").Value = 1470.71432
$ws.Range("N122").Value = -51397

$ws.Range("H132").Value = 1321.5714
$ws.Range("I132").Value = 1269.3846
$ws.Range("K132").Value = 3808.1538
$ws.Range("M132").Value = -1278.1538

$ws.Range("H134").Value = 1526.421
$ws.Range("I134").Value = 1357.7142
$ws.Range("K134").Value = 4073.1426
$ws.Range("M134").Value = -1538.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 297
$ws.Range("J33").Value = 432
$ws.Range("L33").Value = 2592
$ws.Range("N33").Value = -3158

$ws.Range("H42").Value = 3749.5
$ws.Range("J42").Value = 4499
$ws.Range("L42").Value = 13497
$ws.Range("N42").Value = -14565

$ws.Range("H133").Value = 13953.8
$ws.Range("I133").Value = 12182.143
$ws.Range("J133").Value = 14907.77
$ws.Range("K133").Value = 36546.429
$ws.Range("L133").Value = 44723.31
$ws.Range("M133").Value = -31486.429
$ws.Range("N133").Value = -54843.31

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2999.75
$ws.Range("I17").Value = 4499.5
$ws.Range("K17").Value = 4499.5
$ws.Range("M17").Value = -4331.5

$ws.Range("H20").Value = 25000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H26").Value = 19038
$ws.Range("I26").Value = 19038
$ws.Range("K26").Value = 19038
$ws.Range("M26").Value = -18758

$ws.Range("H50").Value = 19038
$ws.Range("I50").Value = 19038
$ws.Range("K50").Value = 19038
$ws.Range("M50").Value = -18540

$ws.Range("H70").Value = 4499.6
$ws.Range("I70").Value = 4549.5
$ws.Range("K70").Value = 4549.5
$ws.Range("M70").Value = -4279.5

$ws.Range("H73").Value = 4499.6
$ws.Range("I73").Value = 4549.5
$ws.Range("K73").Value = 4549.5
$ws.Range("M73").Value = -3613.5

$ws.Range("H102").Value = 1212.7
$ws.Range("I102").Value = 1147.7931
$ws.Range("K102").Value = 1147.7931
$ws.Range("M102").Value = 474.2068999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2482.5
$ws.Range("I82").Value = 2332.3333
$ws.Range("J82").Value = 2632.6667
$ws.Range("K82").Value = 2332.3333
$ws.Range("L82").Value = 2632.6667
$ws.Range("M82").Value = -1971.3333
$ws.Range("N82").Value = -3354.6667

$ws.Range("H85").Value = 2482.5
$ws.Range("I85").Value = 2332.3333
$ws.Range("J85").Value = 2632.6667
$ws.Range("K85").Value = 2332.3333
$ws.Range("L85").Value = 2632.6667
$ws.Range("M85").Value = -1084.3333
$ws.Range("N85").Value = -5128.6667

$ws.Range("H100").Value = 4440
$ws.Range("I100").Value = 4440
$ws.Range("K100").Value = 4440
$ws.Range("M100").Value = -3899

$ws.Range("H122").Value = 2936.75
$ws.Range("I122").Value = 3032
$ws.Range("K122").Value = 9096
$ws.Range("M122").Value = -6646

$ws.Range("H132").Value = 5833.636
$ws.Range("I132").Value = 7747
$ws.Range("J132").Value = 4239.1665
$ws.Range("K132").Value = 23241
$ws.Range("L132").Value = 12717.4995
$ws.Range("M132").Value = -20711
$ws.Range("N132").Value = -17777.4995

$ws.Range("H136").Value = 4734.3076
$ws.Range("I136").Value = 4516.9585
$ws.Range("J136").Value = 7342.5
$ws.Range("K136").Value = 13550.8755
$ws.Range("L136").Value = 22027.5
$ws.Range("M136").Value = -11000.8755
$ws.Range("N136").Value = -27127.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 5292.5713
$ws.Range("J74").Value = 5194.6665
$ws.Range("L74").Value = 5194.6665
$ws.Range("N74").Value = -7066.6665

$ws.Range("H77").Value = 5292.5713
$ws.Range("J77").Value = 5194.6665
$ws.Range("L77").Value = 15583.9995
$ws.Range("N77").Value = -24943.9995

$ws.Range("H104").Value = 11928.5
$ws.Range("J104").Value = 11928.5
$ws.Range("L104").Value = 11928.5
$ws.Range("N104").Value = -18916.5

$ws.Range("H113").Value = 11751
$ws.Range("I113").Value = 11000.333
$ws.Range("K113").Value = 33000.999
$ws.Range("M113").Value = -30830.999

$ws.Range("H122").Value = 3607.75
$ws.Range("J122").Value = 6998.6665
$ws.Range("L122").Value = 20995.9995
$ws.Range("N122").Value = -25895.9995

$ws.Range("H126").Value = 2679.5557
$ws.Range("I126").Value = 2221.3125
$ws.Range("K126").Value = 6663.9375
$ws.Range("M126").Value = -4193.9375
